# Update type / date-of-birth logic
# - Patient sheet: D2 date-of-birth value moves from 1998-02-22 to 2024-02-22
# - Patient sheet: I2/I3 "MRN"-like numeric 123 values become text "G123"
# - Patient tab becomes the active/selected tab (was Guardian)
# - Selection on Guardian moves to C2, selection on Patient moves to K4

$wb = $excel.ActiveWorkbook

$wsGuardian = $wb.Worksheets.Item("Guardian")
$wsPatient  = $wb.Worksheets.Item("Patient")

# Patient!D2 - DateOfBirth column, new serial date 45344 (2024-02-22)
$wsPatient.Range("D2").Value = 45344

# Patient!I2 / I3 - switch from numeric 123 to the text value "G123"
$wsPatient.Range("I2").Value = "G123"
$wsPatient.Range("I3").Value = "G123"

# Update the selections that were saved with the workbook
$wsGuardian.Range("C2").Select()
$wsPatient.Range("K4").Select()

# Patient becomes the active sheet/tab
$wsPatient.Activate()
